$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("openmpi").Name = "Data"
$wb.Worksheets.Item("1cpu").Name = "1core"
$wb.Worksheets.Item("2cpu").Name = "2cores"
$wb.Worksheets.Item("4cpu").Name = "4cores"
$wb.Worksheets.Item("8cpu").Name = "8cores"
